$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: LOTO (SOPs)
$ws.Range("H3").Value = -29
$ws.Range("I3").Formula = "=""04-Nov-2025"""

# Row 4: Endangered by Electricity A safety Training (SOPs)
$ws.Range("H4").Value = -91
$ws.Range("I4").Formula = "=""04-Nov-2025"""

# Row 5: IS0 55001 (Other Trainings)
$ws.Range("H5").Value = 286
$ws.Range("I5").Formula = "=""04-Nov-2025"""

# Convert the helper formulas above into plain literal text values so the
# "LAST UPDATE" cells store "04-Nov-2025" as text (not a live formula),
# matching how the rest of the sheet stores its dates.
$textRange = $ws.Range("I3:I5")
$textRange.Copy()
$textRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0
